# Auto-generated: update FFXIV leve profit computed columns (H-N) across 8 class sheets
# per upstream data refresh (scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 247.59259
$ws.Range("I33").Value = 197.73914
$ws.Range("J33").Value = 534.25
$ws.Range("K33").Value = 197.73914
$ws.Range("L33").Value = 534.25
$ws.Range("M33").Value = 31.26086000000001
$ws.Range("N33").Value = -992.25
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 8137927.5
$ws.Range("I132").Value = 11910690
$ws.Range("J132").Value = 11977.308
$ws.Range("K132").Value = 35732070
$ws.Range("L132").Value = 35931.924
$ws.Range("M132").Value = -35729540
$ws.Range("N132").Value = -40991.924
$ws.Range("H137").Value = 1685.9429
$ws.Range("I137").Value = 1295
$ws.Range("K137").Value = 3885
$ws.Range("M137").Value = -1335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 699.25
$ws.Range("I26").Value = 399
$ws.Range("J26").Value = 1600
$ws.Range("K26").Value = 399
$ws.Range("L26").Value = 1600
$ws.Range("M26").Value = -69
$ws.Range("N26").Value = -2260
$ws.Range("H92").Value = 5000000
$ws.Range("J92").Value = 5000000
$ws.Range("L92").Value = 5000000
$ws.Range("N92").Value = -5004992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1723.2667
$ws.Range("I20").Value = 1445.6957
$ws.Range("J20").Value = 2635.2856
$ws.Range("K20").Value = 1445.6957
$ws.Range("L20").Value = 2635.2856
$ws.Range("M20").Value = -1198.6957
$ws.Range("N20").Value = -3129.2856
$ws.Range("H36").Value = 438.75
$ws.Range("J36").Value = 418
$ws.Range("L36").Value = 418
$ws.Range("N36").Value = -1486
$ws.Range("H75").Value = 4247.3335
$ws.Range("I75").Value = 4247.3335
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 4247.3335
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -3311.3335
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 4247.3335
$ws.Range("I78").Value = 4247.3335
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 12742.0005
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -8062.000499999998
$ws.Range("N78").Value = 0
$ws.Range("H107").Value = 1650.1666
$ws.Range("I107").Value = 1199.5555
$ws.Range("J107").Value = 3002
$ws.Range("K107").Value = 1199.5555
$ws.Range("L107").Value = 3002
$ws.Range("M107").Value = 720.4445000000001
$ws.Range("N107").Value = -6842
$ws.Range("H134").Value = 15374
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 15374
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 46122
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -51192

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1340.2285
$ws.Range("I31").Value = 1167.4517
$ws.Range("J31").Value = 2679.25
$ws.Range("K31").Value = 1167.4517
$ws.Range("L31").Value = 2679.25
$ws.Range("M31").Value = -872.4517000000001
$ws.Range("N31").Value = -3269.25
$ws.Range("H34").Value = 1340.2285
$ws.Range("I34").Value = 1167.4517
$ws.Range("J34").Value = 2679.25
$ws.Range("K34").Value = 1167.4517
$ws.Range("L34").Value = 2679.25
$ws.Range("M34").Value = -965.4517000000001
$ws.Range("N34").Value = -3083.25
$ws.Range("H92").Value = 44300.25
$ws.Range("J92").Value = 44300.25
$ws.Range("L92").Value = 44300.25
$ws.Range("N92").Value = -49292.25
$ws.Range("H94").Value = 1439.091
$ws.Range("I94").Value = 1217.2
$ws.Range("J94").Value = 1624
$ws.Range("K94").Value = 1217.2
$ws.Range("L94").Value = 1624
$ws.Range("M94").Value = -766.2
$ws.Range("N94").Value = -2526
$ws.Range("H99").Value = 1483.1538
$ws.Range("I99").Value = 1503.375
$ws.Range("J99").Value = 1450.8
$ws.Range("K99").Value = 1503.375
$ws.Range("L99").Value = 1450.8
$ws.Range("M99").Value = -5.375
$ws.Range("N99").Value = -4446.8
$ws.Range("H126").Value = 1483.1538
$ws.Range("I126").Value = 1503.375
$ws.Range("J126").Value = 1450.8
$ws.Range("K126").Value = 4510.125
$ws.Range("L126").Value = 4352.4
$ws.Range("M126").Value = -2040.125
$ws.Range("N126").Value = -9292.4
$ws.Range("H134").Value = 22729736
$ws.Range("I134").Value = 2787
$ws.Range("J134").Value = 100001360
$ws.Range("K134").Value = 8361
$ws.Range("L134").Value = 300004080
$ws.Range("M134").Value = -5826
$ws.Range("N134").Value = -300009150
$ws.Range("H141").Value = 872635.7
$ws.Range("J141").Value = 872635.7
$ws.Range("L141").Value = 872635.7
$ws.Range("N141").Value = -882995.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 600
$ws.Range("I15").Value = 600
$ws.Range("K15").Value = 1800
$ws.Range("M15").Value = -1660
$ws.Range("H39").Value = 2963.7778
$ws.Range("J39").Value = 2914.5881
$ws.Range("L39").Value = 8743.764299999999
$ws.Range("N39").Value = -9331.764299999999
$ws.Range("H55").Value = 2376.923
$ws.Range("J55").Value = 2940
$ws.Range("L55").Value = 8820
$ws.Range("N55").Value = -9174
$ws.Range("H97").Value = 980
$ws.Range("I97").Value = 768
$ws.Range("J97").Value = 1245
$ws.Range("K97").Value = 2304
$ws.Range("L97").Value = 3735
$ws.Range("M97").Value = -1808
$ws.Range("N97").Value = -4727
$ws.Range("H107").Value = 12950.375
$ws.Range("I107").Value = 606.6667
$ws.Range("J107").Value = 20356.6
$ws.Range("K107").Value = 1820.0001
$ws.Range("L107").Value = 61069.8
$ws.Range("M107").Value = 99.99990000000003
$ws.Range("N107").Value = -64909.8
$ws.Range("H131").Value = 24394032
$ws.Range("J131").Value = 4165.4053
$ws.Range("L131").Value = 12496.2159
$ws.Range("N131").Value = -22576.2159
$ws.Range("H139").Value = 2033.5
$ws.Range("I139").Value = 2371
$ws.Range("J139").Value = 1599.5714
$ws.Range("K139").Value = 7113
$ws.Range("L139").Value = 4798.7142
$ws.Range("M139").Value = -1973
$ws.Range("N139").Value = -15078.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H122").Value = 1806.9231
$ws.Range("I122").Value = 1464
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 4392
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -1942
$ws.Range("N122").Value = -13750
$ws.Range("H132").Value = 3323.75
$ws.Range("I132").Value = 3254.8125
$ws.Range("J132").Value = 3599.5
$ws.Range("K132").Value = 9764.4375
$ws.Range("L132").Value = 10798.5
$ws.Range("M132").Value = -7234.4375
$ws.Range("N132").Value = -15858.5
$ws.Range("H133").Value = 49735
$ws.Range("J133").Value = 49735
$ws.Range("L133").Value = 49735
$ws.Range("N133").Value = -59855
$ws.Range("H135").Value = 41360
$ws.Range("J135").Value = 32720
$ws.Range("L135").Value = 32720
$ws.Range("N135").Value = -42860
$ws.Range("H137").Value = 50746.668
$ws.Range("J137").Value = 50746.668
$ws.Range("L137").Value = 50746.668
$ws.Range("N137").Value = -60946.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2883.2083
$ws.Range("J40").Value = 5623.8335
$ws.Range("L40").Value = 5623.8335
$ws.Range("N40").Value = -5895.8335
$ws.Range("H68").Value = 1835.2
$ws.Range("I68").Value = 1835.2
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1835.2
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1086.2
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1835.2
$ws.Range("I71").Value = 1835.2
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9176
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -5432
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 103425.09
$ws.Range("I132").Value = 37091
$ws.Range("J132").Value = 128300.375
$ws.Range("K132").Value = 111273
$ws.Range("L132").Value = 384901.125
$ws.Range("M132").Value = -108743
$ws.Range("N132").Value = -389961.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 484.3
$ws.Range("I113").Value = 329.2
$ws.Range("J113").Value = 639.4
$ws.Range("K113").Value = 987.5999999999999
$ws.Range("L113").Value = 1918.2
$ws.Range("M113").Value = 1182.4
$ws.Range("N113").Value = -6258.2
